# Refresh the scraped "想去人数" (want-to-go count) figures in column F
# across all four sheets, matching the data pulled at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (exhibitions)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1402
$ws.Range("F4").Value = 25696
$ws.Range("F5").Value = 566
$ws.Range("F7").Value = 569
$ws.Range("F8").Value = 161
$ws.Range("F12").Value = 196
$ws.Range("F13").Value = 170
$ws.Range("F14").Value = 42
$ws.Range("F15").Value = 272
$ws.Range("F16").Value = 329
$ws.Range("F18").Value = 1463
$ws.Range("F19").Value = 152

# Sheet 2: 演出 (performances)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 216
$ws.Range("F6").Value = 67
$ws.Range("F10").Value = 414
$ws.Range("F14").Value = 16
$ws.Range("F15").Value = 22

# Sheet 3: 本地生活 (local life)
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 4875
$ws.Range("F3").Value = 170

# Sheet 4: 全部类型 (all types - union of the above three sheets)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1402
$ws.Range("F4").Value = 4875
$ws.Range("F5").Value = 170
$ws.Range("F6").Value = 25696
$ws.Range("F7").Value = 566
$ws.Range("F10").Value = 216
$ws.Range("F11").Value = 569
$ws.Range("F14").Value = 161
$ws.Range("F15").Value = 67
$ws.Range("F16").Value = 67
$ws.Range("F20").Value = 414
$ws.Range("F25").Value = 196
$ws.Range("F26").Value = 170
$ws.Range("F27").Value = 42
$ws.Range("F29").Value = 272
$ws.Range("F31").Value = 16
$ws.Range("F32").Value = 329
$ws.Range("F34").Value = 22
$ws.Range("F35").Value = 1463
$ws.Range("F36").Value = 152
